$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet (and the matching defined name "VANS_leverandører") are named
# after the "last updated" date shown in the tab. This edit bumps that
# date from 02-12-2025 to 05-12-2025. Renaming the worksheet also updates
# the workbook-scoped defined name that references it.
$ws.Name = "Opdateret d. 05-12-2025"
